$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new TOTAL / ASOPAGOS headers -----------------------------------
$totalAso = $ws.Range("O1:P1")
$totalAso.Font.Bold = $true
$totalAso.Font.Name = "Arial"
$totalAso.Font.Size = 12
$totalAso.WrapText = $true
$totalAso.Borders.LineStyle = 1
$ws.Range("P1").Value = "ASOPAGOS"
$ws.Range("O1").Value = "TOTAL"

# --- Row 1: new ENTIDAD header spanning A1:C1 ------------------------------
$entidad = $ws.Range("A1:C1")
$entidad.Font.Bold = $true
$entidad.Font.Name = "Arial"
$entidad.Font.Size = 12
$entidad.HorizontalAlignment = -4108
$entidad.WrapText = $true
$entidad.Borders.LineStyle = 1
$entidad.Merge()
$ws.Range("A1").Value = "ENTIDAD"

# --- Row 2: rename UNIDAD 02/08/09 -> UNIDAD 2/8 /9 -----------------------
$ws.Range("D2").Value = "UNIDAD 2"
$ws.Range("E2").Value = "UNIDAD 8 "
$ws.Range("F2").Value = "UNIDAD 9"

# --- sheet view tidy up -----------------------------------------------------
$ws.Range("D12").Select()
